$wb = $excel.ActiveWorkbook

# The "Italy" sheet is the template for every country tab in this workbook.
# Spain's test data is added the same way the other markets were: duplicate
# the Italy tab, rename it, and swap in the Spain-specific values.
$italy = $wb.Worksheets.Item("Italy")

# Make sure Italy's own selection is reset to the full used range before we
# branch off the copy (this matches how the sheet settles once it stops
# being the active tab).
$italy.Activate() | Out-Null
$italy.Range("A1:D14").Select() | Out-Null

# Duplicate "Italy" and place the copy right after it.
$italy.Copy([System.Reflection.Missing]::Value, $italy) | Out-Null

$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# Market name + printer id for the Spain tab.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2051 "

# Make Spain the active tab with B4:B5 selected.
$spain.Activate() | Out-Null
$spain.Range("B4:B5").Select() | Out-Null
